$d = $word.ActiveDocument

# wdReplaceAll = 2 ; wdFindContinue = 1 (we scope searches to a single
# paragraph's Range so "replace all" only ever touches that paragraph).

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $r = $d.Paragraphs($paraIndex).Range
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "MISS para=$paraIndex find=$findText"
    }
}

# --- Requirements paragraph (originally paragraph 10) ---------------------
Replace-InParagraph 10 " classes: " " classes and functions: "
Replace-InParagraph 10 "Next, it will require a class for player selection (chara_select) that will set what characters player 1 and player 2 are. Third, " "Next, "
Replace-InParagraph 10 "typical anime shounen " "typical anime Shonen "

# --- Main() heading (originally paragraph 11) ------------------------------
Replace-InParagraph 11 "Main():" "Main() (This is a function):"

# --- "2. Wiki" menu heading (originally paragraph 20) ---------------------
Replace-InParagraph 20 "2. Wiki" "2. Character Wiki"

# --- Remove the Chara_select() heading + description paragraphs -----------
# (originally paragraphs 24 and 25); delete from the bottom up so indices
# of earlier paragraphs are unaffected.
$d.Paragraphs(25).Range.Delete()
$d.Paragraphs(24).Range.Delete()

# --- Insert two blank paragraphs after the Wiki() description -------------
# After the deletions above, "Contains information on how to use each
# character. Returns strings. " is paragraph 27 (was 29).
$wikiDescRange = $d.Paragraphs(27).Range
$wikiDescRange.InsertParagraphAfter()
$d.Paragraphs(28).Range.InsertParagraphAfter()

# --- Samurai() description text (now paragraph 33) ------------------------
Replace-InParagraph 33 "unless is blocked), block (starts with 5, -1 each use, protects from attack). With slash, if enemy also attacks at same time" "unless is blocked. If it is blocked, then the sword breaks and unsheathed becomes False), block (starts with 5, -1 each use, protects from attack). With slash, if enemy also attacks at same time"

Write-Output "done"
